# Adiciona medidas de rotação e translação a dados
# Adds rotation/translation helper columns (J:O) to "Planetas e Sol",
# a new moon data point in row 11, an extra computed row for row 12,
# and re-applies the border-carrying style to the "ring info" caption rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planetas e Sol")

# ---- Period (I) -> log / exp helper columns, and new "L" rotation data
#      with its own log / exp helper columns (M:O), for rows 3..10 ----

$rotation = @{
    3  = 58
    4  = 243.025
    5  = 1
    6  = 1
    7  = 0.333
    8  = 0.4
    9  = 0.71
    10 = 0.67
}

foreach ($r in 3..10) {
    $ws.Range("J$r").Formula = "=ROUND(LOG(I$r,10), 1)"
    $ws.Range("K$r").Formula = "=ROUND(EXP(J$r),2)"
    $ws.Range("L$r").Value = $rotation[$r]
    $ws.Range("M$r").Formula = "=L$r*10"
    $ws.Range("N$r").Formula = "=LOG10(M$r)"
    $ws.Range("O$r").Formula = "=ROUND(EXP(N$r),2)"
}

# ---- New row 11: an extra moon-style entry with just I/J/K populated ----
$ws.Range("I11").Value = 28
$ws.Range("J11").Formula = "=ROUND(LOG(I11,10), 1)"
$ws.Range("K11").Formula = "=ROUND(EXP(J11),2)"

# ---- Row 12 gains an I/J/K computation (human lifetime in days) ----
$ws.Range("I12").Formula = "=75*365"
$ws.Range("J12").Formula = "=ROUND(LOG(I12,10), 1)"
$ws.Range("K12").Formula = "=ROUND(EXP(J12),2)"

# ---- Re-stamp the caption rows (merged B12:D12 / B13:D13) with an
#      explicit (empty) border so they pick up the bordered style ----
$ws.Range("B12:D13").Borders.LineStyle = 1
$ws.Range("B12:D13").Borders.LineStyle = 0

$ws.Range("I13").Select()
